# Refresh the legacy GSC export "Chart" sheet: the rolling date window
# advanced by one day, so the oldest row (2025-10-30) drops off the top
# and every later row's data shifts up to fill in, just like a real
# re-export of the report would look once re-pasted into the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Deleting the entire row shifts everything below it up by one and lets
# Excel renumber/reshuffle the shared strings & dimension automatically.
$ws.Range("A2").EntireRow.Delete()
